# Add a new patient row (row 3) to the "Bemorlar ro`yxati" (Patients list) sheet.
# Columns: A=ID, B=Name, C=Phone, D=Floor, E=Room, F=Registration date, G=Status
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the values that would otherwise be auto-detected as a
# number ("6"), a phone number ("+998909999999") or a date ("2022-12-07") so
# that they are written out exactly as typed (same as the other text columns).
$ws.Range("A3").NumberFormat = "@"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("F3").NumberFormat = "@"

$ws.Range("A3").Value = "6"
$ws.Range("B3").Value = "Asal"
$ws.Range("C3").Value = "+998909999999"
$ws.Range("D3").Value = "1qavat"
$ws.Range("E3").Value = "21xona"
$ws.Range("F3").Value = "2022-12-07"
$ws.Range("G3").Value = "ACTIVE"

# Drop the temporary text number-format again so the new row keeps the same
# (default) cell style as the rest of the sheet.
$ws.Range("A3:G3").ClearFormats()
